# Feedback 2 General is done
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 22 with the new "Feedback 2" entry
$ws.Range("A22").Value = "App icon + Rotate"
$ws.Range("B22").Value = (Get-Date -Year 2018 -Month 1 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0.25"
$ws.Range("C22").NumberFormat = "0"
$ws.Range("D22").Value = "I improved the app icon. Now it hasn't a black border anymore. And a solved the rotating issue."

# Update the view/selection state to match the author's last position
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("A23").Select()
